# Update cryptocurrency price/volume data and re-rank coin list
# (author commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.179.12"
$ws.Range("E2").Value = "  -0.38%  "

# Row 3
$ws.Range("D3").Value = "1.910.77"
$ws.Range("E3").Value = "  +1.98%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'314.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.93%  "

# Row 6
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("D7").Value = "'0.5093"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.79%  "

# Row 8
$ws.Range("D8").Value = "'0.3922"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.32%  "

# Row 9
$ws.Range("D9").Value = "'0.09249"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.86%  "

# Row 10
$ws.Range("D10").Value = "'1.138"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.53%  "

# Row 11
$ws.Range("D11").Value = "'41.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.46%  "

# Row 12
$ws.Range("E12").Value = "  -1.55%  "

# Row 13
$ws.Range("D13").Value = "'20.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.74%  "

# Row 14
$ws.Range("D14").Value = "1.906.92"
$ws.Range("E14").Value = "  +1.62%  "

# Row 15
$ws.Range("D15").Value = "'7.317"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.57%  "

# Row 16
$ws.Range("D16").Value = "'1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "

# Row 17
$ws.Range("D17").Value = "'0.00001119"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.21%  "

# Row 18
$ws.Range("D18").Value = "'92.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.69%  "

# Row 19
$ws.Range("D19").Value = "'0.06611"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "

# Row 20
$ws.Range("D20").Value = "'17.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.76%  "

# Row 21
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'6.230"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.70%  "

# Row 23
$ws.Range("D23").Value = "28.234.72"
$ws.Range("E23").Value = "  -0.36%  "

# Row 24
$ws.Range("E24").Value = "  +1.63%  "

# Row 25
$ws.Range("E25").Value = "  +1.62%  "

# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.590"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.91%  "

# Row 27
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.125.21"
$ws.Range("E27").Value = "  +1.47%  "

# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'21.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.87%  "

# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'158.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.38%  "

# Row 30
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'127.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.48%  "

# Row 31
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.095"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.54%  "

# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.1076"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.76%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.637"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.616"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.14%  "

# Row 35
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "'9.693"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.01%  "

# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.06666"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.01%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02427"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.05%  "

# Row 38
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.240"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.76%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2195"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "

# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.287"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.56%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6461"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.58%  "

# Row 42
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'11.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.19%  "

# Row 43
$ws.Range("B43").Value = "InternetComputer(DFINITY)"
$ws.Range("C43").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D43").Value = "'4.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.05%  "

# Row 44
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.08%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6058"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.05%  "

# Row 47
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.720"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.67%  "

# Row 48
$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'1.287"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.57%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'2.013"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.34%  "

# Row 50
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'123.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.63%  "

# Row 51
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.187"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.76%  "
